$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 1.03863614
$ws.Range("D3").Value = 2.3255502
